$d = $word.ActiveDocument

# 1. "Current file not working" -> "Good"
$d.Content.Find.Execute("Current file not working", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Good", 2)

# 2. "Works (but I don't think its working correctly)" -> "Good"
$d.Content.Find.Execute("Works (but I don't think its working correctly)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Good", 2)

# 3. "Same as test4 test5" -> "Fails"
$d.Content.Find.Execute("Same as test4 test5", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fails", 2)
